{"js": "// The document contains a single table of two-digit \u00f7 one-digit division\n// practice problems. Every populated cell's answer text is replaced with a\n// new problem/answer string. Because a couple of the old strings repeat\n// (e.g. \"43\u00f74=10, 3\" appears twice but maps to two different replacements),\n// cells are addressed by their (row, column) position in the table rather\n// than by searching for the old text.\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Map of row index -> array of 5 new cell values (row indices are 0-based\n// and count every table row, including the blank spacer rows).\nconst newValues = {\n  0: [\"37\u00f76=6, 1\", \"67\u00f79=7, 4\", \"10\u00f75=2, 0\", \"60\u00f76=10, 0\", \"66\u00f72=33, 0\"],\n  4: [\"30\u00f76=5, 0\", \"88\u00f76=14, 4\", \"69\u00f77=9, 6\", \"15\u00f79=1, 6\", \"75\u00f74=18, 3\"],\n  8: [\"98\u00f73=32, 2\", \"24\u00f73=8, 0\", \"79\u00f73=26, 1\", \"81\u00f79=9, 0\", \"43\u00f72=21, 1\"],\n  12: [\"99\u00f72=49, 1\", \"51\u00f74=12, 3\", \"77\u00f72=38, 1\", \"43\u00f77=6, 1\", \"93\u00f79=10, 3\"],\n  16: [\"75\u00f77=10, 5\", \"97\u00f77=13, 6\", \"24\u00f78=3, 0\", \"83\u00f72=41, 1\", \"25\u00f73=8, 1\"],\n};\n\nfor (const [rowStr, values] of Object.entries(newValues)) {\n  const row = Number(rowStr);\n  for (let col = 0; col < values.length; col++) {\n    table.getCell(row, col).value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single table of two-digit \u00f7 one-digit division\n# practice problems. Every populated cell's answer text is replaced with a\n# new problem/answer string. Because a couple of the old strings repeat\n# (e.g. \"43\u00f74=10, 3\" appears twice but maps to two different replacements),\n# cells are addressed by their (row, column) position in the table rather\n# than by searching for the old text. Word COM table indices are 1-based.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"37\u00f76=6, 1\", \"67\u00f79=7, 4\", \"10\u00f75=2, 0\", \"60\u00f76=10, 0\", \"66\u00f72=33, 0\")\n    5  = @(\"30\u00f76=5, 0\", \"88\u00f76=14, 4\", \"69\u00f77=9, 6\", \"15\u00f79=1, 6\", \"75\u00f74=18, 3\")\n    9  = @(\"98\u00f73=32, 2\", \"24\u00f73=8, 0\", \"79\u00f73=26, 1\", \"81\u00f79=9, 0\", \"43\u00f72=21, 1\")\n    13 = @(\"99\u00f72=49, 1\", \"51\u00f74=12, 3\", \"77\u00f72=38, 1\", \"43\u00f77=6, 1\", \"93\u00f79=10, 3\")\n    17 = @(\"75\u00f77=10, 5\", \"97\u00f77=13, 6\", \"24\u00f78=3, 0\", \"83\u00f72=41, 1\", \"25\u00f73=8, 1\")\n}\n\nforeach ($row in $newValues.Keys) {\n    $vals = $newValues[$row]\n    for ($c = 0; $c -lt $vals.Count; $c++) {\n        $t.Cell($row, $c + 1).Range.Text = $vals[$c]\n    }\n}\n"}
